$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A2").Value  = "Pictures_Practice/PICTURE_612.png"
$ws.Range("A3").Value  = "Pictures_Practice/PICTURE_599.png"
$ws.Range("A4").Value  = "Pictures_Practice/PICTURE_570.png"
$ws.Range("A5").Value  = "Pictures_Practice/PICTURE_570.png"
$ws.Range("A6").Value  = "Pictures_Practice/PICTURE_733.png"
$ws.Range("A7").Value  = "Pictures_Practice/PICTURE_614.png"
$ws.Range("A8").Value  = "Pictures_Practice/PICTURE_110.png"
$ws.Range("A9").Value  = "Pictures_Practice/PICTURE_110.png"
$ws.Range("A10").Value = "Pictures_Practice/PICTURE_599.png"
$ws.Range("A11").Value = "Pictures_Practice/PICTURE_12.png"
$ws.Range("A12").Value = "Pictures_Practice/PICTURE_12.png"
$ws.Range("A13").Value = "Pictures_Practice/PICTURE_612.png"
